$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 3-16 (columns A: Player, B: Position, C: Team)
$data = @(
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Jusuf Nurkic", "C", "Phoenix Suns"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Moritz Wagner", "C", "Orlando Magic"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Brook Lopez", "C", "Milwaukee Bucks")
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
